$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Skill")

# Remove the old animset/attack_type/cooltime/abnormal/type columns (D:H) and
# replace them with the new style_tube / enhancer_tube / cooler_tube columns.
$rng = $ws.Range($ws.Cells.Item(1,4), $ws.Cells.Item(7,8))
$rng.EntireColumn.Delete()

# Header row
$ws.Range("D1").Value = "style_tube"
$ws.Range("E1").Value = "enhancer_tube"
$ws.Range("F1").Value = "cooler_tube"

# 19101 - DeadlyAttack
$ws.Range("D2").Value = "no107_style"
$ws.Range("E2").Value = "no107_enhancer"
$ws.Range("F2").Value = "no107_cooler"

# 19102 - CriticalHit
$ws.Range("D3").Value = "no108_style"
$ws.Range("E3").Value = "no108_enhancer"
$ws.Range("F3").Value = "no108_cooler"

# 19103 - Badbomb
$ws.Range("D4").Value = "jake_style"
$ws.Range("E4").Value = "jake_enhancer"
$ws.Range("F4").Value = "jake_cooler"

# Rows 5-7 (19200/19201/19202) have no style/enhancer/cooler tube values.

$ws.Columns.Item(4).ColumnWidth = 10.65
$ws.Columns.Item(5).ColumnWidth = 15.25
$ws.Columns.Item(6).ColumnWidth = 12.25

[void]$ws.Range("B3").Select()
